$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 339
$ws1.Range("F4").Value = 1418
$ws1.Range("F5").Value = 678

# Sheet "全部类型" (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 339
$ws4.Range("F4").Value = 1418
$ws4.Range("F6").Value = 678
